$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the "Periodo Mora" data table (rows 16-53) in ascending period order,
# interleaving NELCY LINEY BUELVAS URANGO and RAUL DARIO BUELVAS URANGO for each period.
$data = @(
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1612', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1612', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1701', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1701', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1702', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1702', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1703', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1703', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1704', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1704', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1705', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1705', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1706', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1706', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1707', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1707', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1708', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1708', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1709', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1709', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1710', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1710', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1711', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1711', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1712', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1712', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1801', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1801', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1802', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1802', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1803', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1803', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1804', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1804', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1805', 27578, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1805', 27578, 689455)
    ,@('CC', '45535199', 'NELCY LINEY BUELVAS URANGO', '1806', 16547, 689455)
    ,@('CC', '73206843', 'RAUL DARIO BUELVAS URANGO', '1806', 16547, 689455)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
